$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02743666666666666
$ws.Range("H2").Value = 0.08231
$ws.Range("I2").Value = 0.007366285056527356
$ws.Range("J2").Value = 0.007366285056527356
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 0.312815
$ws.Range("N2").Value = 0.938445
$ws.Range("O2").Value = 0.0082131704949067
$ws.Range("P2").Value = 0.0082131704949067
$ws.Range("Q2").Value = 0.008582600883333333
$ws.Range("R2").Value = 0.07724340795
$ws.Range("S2").Value = 0.00006050055508334261
$ws.Range("T2").Value = 0.00006050055508334261

$ws.Range("G3").Value = 0.02743666666666666
$ws.Range("H3").Value = 0.08231
$ws.Range("I3").Value = 0.007366285056527356
$ws.Range("J3").Value = 0.007366285056527356
$ws.Range("O3").Value = 0.6524076620340182
$ws.Range("P3").Value = 0.6524076620340182
$ws.Range("Q3").Value = 0.6817531159177778
$ws.Range("R3").Value = 6.13577804326
$ws.Range("S3").Value = 0.004805820811605138
$ws.Range("T3").Value = 0.004805820811605138

$ws.Range("G4").Value = 0.02743666666666666
$ws.Range("H4").Value = 0.08231
$ws.Range("I4").Value = 0.007366285056527356
$ws.Range("J4").Value = 0.007366285056527356
$ws.Range("M4").Value = 12.866992
$ws.Range("N4").Value = 38.600976
$ws.Range("O4").Value = 0.3378316226926476
$ws.Range("P4").Value = 0.3378316226926476
$ws.Range("Q4").Value = 0.3530273705066667
$ws.Range("R4").Value = 3.17724633456
$ws.Range("S4").Value = 0.002488564033863238
$ws.Range("T4").Value = 0.002488564033863238

$ws.Range("G5").Value = 0.02743666666666666
$ws.Range("H5").Value = 0.08231
$ws.Range("I5").Value = 0.007366285056527356
$ws.Range("J5").Value = 0.007366285056527356
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 0.05894133333333334
$ws.Range("N5").Value = 0.176824
$ws.Range("O5").Value = 0.001547544778427486
$ws.Range("P5").Value = 0.001547544778427486
$ws.Range("Q5").Value = 0.001617153715555556
$ws.Range("R5").Value = 0.01455438344
$ws.Range("S5").Value = 0.00001139965597563733
$ws.Range("T5").Value = 0.00001139965597563733

$ws.Range("G6").Value = 3.368329
$ws.Range("I6").Value = 0.9043398704228307
$ws.Range("J6").Value = 0.9043398704228307
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 0.312815
$ws.Range("N6").Value = 0.938445
$ws.Range("O6").Value = 0.0082131704949067
$ws.Range("P6").Value = 0.0082131704949067
$ws.Range("Q6").Value = 1.053663836135
$ws.Range("R6").Value = 9.482974525214999
$ws.Range("S6").Value = 0.007427497541124541
$ws.Range("T6").Value = 0.007427497541124541

$ws.Range("G7").Value = 3.368329
$ws.Range("I7").Value = 0.9043398704228307
$ws.Range("J7").Value = 0.9043398704228307
$ws.Range("O7").Value = 0.6524076620340182
$ws.Range("P7").Value = 0.6524076620340182
$ws.Range("Q7").Value = 83.69707658314466
$ws.Range("R7").Value = 753.273689248302
$ws.Range("S7").Value = 0.589998260546706
$ws.Range("T7").Value = 0.589998260546706

$ws.Range("G8").Value = 3.368329
$ws.Range("I8").Value = 0.9043398704228307
$ws.Range("J8").Value = 0.9043398704228307
$ws.Range("M8").Value = 12.866992
$ws.Range("N8").Value = 38.600976
$ws.Range("O8").Value = 0.3378316226926476
$ws.Range("P8").Value = 0.3378316226926476
$ws.Range("Q8").Value = 43.340262296368
$ws.Range("R8").Value = 390.062360667312
$ws.Range("S8").Value = 0.3055146058906036
$ws.Range("T8").Value = 0.3055146058906036

$ws.Range("G9").Value = 3.368329
$ws.Range("I9").Value = 0.9043398704228307
$ws.Range("J9").Value = 0.9043398704228307
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 0.05894133333333334
$ws.Range("N9").Value = 0.176824
$ws.Range("O9").Value = 0.001547544778427486
$ws.Range("P9").Value = 0.001547544778427486
$ws.Range("Q9").Value = 0.1985338023653333
$ws.Range("R9").Value = 1.786804221288
$ws.Range("S9").Value = 0.001399506444396641
$ws.Range("T9").Value = 0.001399506444396641

$ws.Range("G10").Value = 0.3288616666666667
$ws.Range("H10").Value = 0.9865849999999999
$ws.Range("I10").Value = 0.08829384452064198
$ws.Range("J10").Value = 0.08829384452064198
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 0.312815
$ws.Range("N10").Value = 0.938445
$ws.Range("O10").Value = 0.0082131704949067
$ws.Range("P10").Value = 0.0082131704949067
$ws.Range("Q10").Value = 0.1028728622583333
$ws.Range("R10").Value = 0.9258557603249999
$ws.Range("S10").Value = 0.0007251723986988163
$ws.Range("T10").Value = 0.0007251723986988163

$ws.Range("G11").Value = 0.3288616666666667
$ws.Range("H11").Value = 0.9865849999999999
$ws.Range("I11").Value = 0.08829384452064198
$ws.Range("J11").Value = 0.08829384452064198
$ws.Range("O11").Value = 0.6524076620340182
$ws.Range("P11").Value = 0.6524076620340182
$ws.Range("Q11").Value = 8.171636470267778
$ws.Range("R11").Value = 73.54472823241
$ws.Range("S11").Value = 0.05760358067570715
$ws.Range("T11").Value = 0.05760358067570715

$ws.Range("G12").Value = 0.3288616666666667
$ws.Range("H12").Value = 0.9865849999999999
$ws.Range("I12").Value = 0.08829384452064198
$ws.Range("J12").Value = 0.08829384452064198
$ws.Range("M12").Value = 12.866992
$ws.Range("N12").Value = 38.600976
$ws.Range("O12").Value = 0.3378316226926476
$ws.Range("P12").Value = 0.3378316226926476
$ws.Range("Q12").Value = 4.231460434106667
$ws.Range("R12").Value = 38.08314390696
$ws.Range("S12").Value = 0.02982845276818081
$ws.Range("T12").Value = 0.02982845276818081

$ws.Range("G13").Value = 0.3288616666666667
$ws.Range("H13").Value = 0.9865849999999999
$ws.Range("I13").Value = 0.08829384452064198
$ws.Range("J13").Value = 0.08829384452064198
$ws.Range("K13").Value = 3.0
$ws.Range("L13").Value = 1.0
$ws.Range("M13").Value = 0.05894133333333334
$ws.Range("N13").Value = 0.176824
$ws.Range("O13").Value = 0.001547544778427486
$ws.Range("P13").Value = 0.001547544778427486
$ws.Range("Q13").Value = 0.01938354511555556
$ws.Range("R13").Value = 0.17445190604
$ws.Range("S13").Value = 0.0001366386780552078
$ws.Range("T13").Value = 0.0001366386780552078
